# Auto-generated edit script: updates market-price / profit columns (H:N)
# across multiple job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# a scheduled data-refresh snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1772.8667
$ws.Range("I29").Value = 198.25
$ws.Range("J29").Value = 2345.4546
$ws.Range("K29").Value = 594.75
$ws.Range("L29").Value = 7036.3638
$ws.Range("M29").Value = -313.75
$ws.Range("N29").Value = -7598.3638
$ws.Range("H38").Value = 4152
$ws.Range("J38").Value = 5328.5713
$ws.Range("L38").Value = 15985.7139
$ws.Range("N38").Value = -16729.7139
$ws.Range("H100").Value = 2296.8276
$ws.Range("I100").Value = 2311.8462
$ws.Range("J100").Value = 2166.6667
$ws.Range("K100").Value = 2311.8462
$ws.Range("L100").Value = 2166.6667
$ws.Range("M100").Value = -1770.8462
$ws.Range("N100").Value = -3248.6667
$ws.Range("H132").Value = 9266811
$ws.Range("I132").Value = 11500743
$ws.Range("K132").Value = 34502229
$ws.Range("M132").Value = -34499699
$ws.Range("H137").Value = 1084.1
$ws.Range("I137").Value = 698.38464
$ws.Range("J137").Value = 1800.4286
$ws.Range("K137").Value = 2095.15392
$ws.Range("L137").Value = 5401.2858
$ws.Range("M137").Value = 454.8460800000003
$ws.Range("N137").Value = -10501.2858
$ws.Range("H138").Value = 1144.4166
$ws.Range("I138").Value = 625.46344
$ws.Range("J138").Value = 1531.2727
$ws.Range("K138").Value = 1876.39032
$ws.Range("L138").Value = 4593.8181
$ws.Range("M138").Value = 3263.60968
$ws.Range("N138").Value = -14873.8181
$ws.Range("H141").Value = 630.8182
$ws.Range("I141").Value = 630.8182
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1892.4546
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3287.5454
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3560.1956
$ws.Range("I32").Value = 3107.0952
$ws.Range("K32").Value = 3107.0952
$ws.Range("M32").Value = -2820.0952
$ws.Range("H61").Value = 26316788
$ws.Range("I61").Value = 40000836
$ws.Range("J61").Value = 1308.6154
$ws.Range("K61").Value = 40000836
$ws.Range("L61").Value = 1308.6154
$ws.Range("M61").Value = -40000624
$ws.Range("N61").Value = -1732.6154
$ws.Range("H132").Value = 1182.4
$ws.Range("I132").Value = 1129.2
$ws.Range("K132").Value = 3387.6
$ws.Range("M132").Value = -857.6000000000004
$ws.Range("H136").Value = 26316788
$ws.Range("I136").Value = 40000836
$ws.Range("J136").Value = 1308.6154
$ws.Range("K136").Value = 120002508
$ws.Range("L136").Value = 3925.8462
$ws.Range("M136").Value = -119999958
$ws.Range("N136").Value = -9025.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1700.0646
$ws.Range("I20").Value = 1308.45
$ws.Range("J20").Value = 2412.0908
$ws.Range("K20").Value = 1308.45
$ws.Range("L20").Value = 2412.0908
$ws.Range("M20").Value = -1061.45
$ws.Range("N20").Value = -2906.0908
$ws.Range("H134").Value = 4038.1794
$ws.Range("I134").Value = 1069.1389
$ws.Range("K134").Value = 3207.4167
$ws.Range("M134").Value = -672.4166999999998
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1512.5892
$ws.Range("I31").Value = 1381.4706
$ws.Range("J31").Value = 2850
$ws.Range("K31").Value = 1381.4706
$ws.Range("L31").Value = 2850
$ws.Range("M31").Value = -1086.4706
$ws.Range("N31").Value = -3440
$ws.Range("H34").Value = 1512.5892
$ws.Range("I34").Value = 1381.4706
$ws.Range("J34").Value = 2850
$ws.Range("K34").Value = 1381.4706
$ws.Range("L34").Value = 2850
$ws.Range("M34").Value = -1179.4706
$ws.Range("N34").Value = -3254
$ws.Range("H58").Value = 700.28815
$ws.Range("I58").Value = 611.0784
$ws.Range("J58").Value = 1269
$ws.Range("K58").Value = 611.0784
$ws.Range("L58").Value = 1269
$ws.Range("M58").Value = -408.0784
$ws.Range("N58").Value = -1675
$ws.Range("H86").Value = 1765049.5
$ws.Range("I86").Value = 2672892.8
$ws.Range("J86").Value = 19197.076
$ws.Range("K86").Value = 2672892.8
$ws.Range("L86").Value = 19197.076
$ws.Range("M86").Value = -2671769.8
$ws.Range("N86").Value = -21443.076
$ws.Range("H89").Value = 1765049.5
$ws.Range("I89").Value = 2672892.8
$ws.Range("J89").Value = 19197.076
$ws.Range("K89").Value = 13364464
$ws.Range("L89").Value = 95985.38
$ws.Range("M89").Value = -13358848
$ws.Range("N89").Value = -107217.38
$ws.Range("H132").Value = 3528.66
$ws.Range("I132").Value = 3823.372
$ws.Range("J132").Value = 1718.2858
$ws.Range("K132").Value = 11470.116
$ws.Range("L132").Value = 5154.857400000001
$ws.Range("M132").Value = -8940.116
$ws.Range("N132").Value = -10214.8574
$ws.Range("H134").Value = 11906154
$ws.Range("I134").Value = 1495.4412
$ws.Range("J134").Value = 62500950
$ws.Range("K134").Value = 4486.3236
$ws.Range("L134").Value = 187502850
$ws.Range("M134").Value = -1951.3236
$ws.Range("N134").Value = -187507920
$ws.Range("H136").Value = 700.28815
$ws.Range("I136").Value = 611.0784
$ws.Range("J136").Value = 1269
$ws.Range("K136").Value = 1833.2352
$ws.Range("L136").Value = 3807
$ws.Range("M136").Value = 716.7647999999999
$ws.Range("N136").Value = -8907
$ws.Range("H3").Value = 692666.3
$ws.Range("I3").Value = 97999
$ws.Range("J3").Value = 990000
$ws.Range("K3").Value = 97999
$ws.Range("L3").Value = 990000
$ws.Range("M3").Value = -97886
$ws.Range("N3").Value = -990226
$ws.Range("H99").Value = 2980
$ws.Range("I99").Value = 2980
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2980
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1482
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 2980
$ws.Range("I126").Value = 2980
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8940
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6470
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1090.3
$ws.Range("I5").Value = 1122.2
$ws.Range("J5").Value = 930.8
$ws.Range("K5").Value = 3366.6
$ws.Range("L5").Value = 2792.4
$ws.Range("M5").Value = -3254.6
$ws.Range("N5").Value = -3016.4
$ws.Range("H122").Value = 779.9032
$ws.Range("I122").Value = 676.7692
$ws.Range("K122").Value = 6090.922799999999
$ws.Range("M122").Value = -3640.922799999999
$ws.Range("H125").Value = 4813.2856
$ws.Range("H130").Value = 1929.8182
$ws.Range("J130").Value = 2019.8
$ws.Range("L130").Value = 6059.4
$ws.Range("N130").Value = -16099.4
$ws.Range("H135").Value = 1090.3
$ws.Range("I135").Value = 1122.2
$ws.Range("J135").Value = 930.8
$ws.Range("K135").Value = 10099.8
$ws.Range("L135").Value = 8377.199999999999
$ws.Range("M135").Value = -7564.800000000001
$ws.Range("N135").Value = -13447.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 28320.334
$ws.Range("J86").Value = 28320.334
$ws.Range("L86").Value = 28320.334
$ws.Range("N86").Value = -30692.334
$ws.Range("H89").Value = 28320.334
$ws.Range("J89").Value = 28320.334
$ws.Range("L89").Value = 84961.00199999999
$ws.Range("N89").Value = -96817.00199999999
$ws.Range("H97").Value = 448.41177
$ws.Range("I97").Value = 479.76923
$ws.Range("J97").Value = 346.5
$ws.Range("K97").Value = 479.76923
$ws.Range("L97").Value = 346.5
$ws.Range("M97").Value = 16.23077000000001
$ws.Range("N97").Value = -1338.5
$ws.Range("H132").Value = 1647.7188
$ws.Range("I132").Value = 1360.2593
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 4080.7779
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -1550.7779
$ws.Range("N132").Value = -14660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2045.1111
$ws.Range("I7").Value = 1915.2858
$ws.Range("K7").Value = 1915.2858
$ws.Range("M7").Value = -1803.2858
$ws.Range("H126").Value = 2045.1111
$ws.Range("I126").Value = 1915.2858
$ws.Range("K126").Value = 5745.857400000001
$ws.Range("M126").Value = -3275.857400000001
$ws.Range("H132").Value = 23542.61
$ws.Range("I132").Value = 1531.7037
$ws.Range("K132").Value = 4595.1111
$ws.Range("M132").Value = -2065.1111
$ws.Range("H136").Value = 3877.7632
$ws.Range("I136").Value = 3953
$ws.Range("K136").Value = 11859
$ws.Range("M136").Value = -9309
$ws.Range("H40").Value = 2899.8
$ws.Range("I40").Value = 2666.4443
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2666.4443
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2530.4443
$ws.Range("N40").Value = -5272

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2879
$ws.Range("I96").Value = 1900.75
$ws.Range("J96").Value = 4183.3335
$ws.Range("K96").Value = 1900.75
$ws.Range("L96").Value = 4183.3335
$ws.Range("M96").Value = -527.75
$ws.Range("N96").Value = -6929.3335
$ws.Range("H126").Value = 55557336
$ws.Range("I126").Value = 71430090
$ws.Range("K126").Value = 214290270
$ws.Range("M126").Value = -214287800
$ws.Range("H132").Value = 1091.4348
$ws.Range("I132").Value = 1121.88
$ws.Range("J132").Value = 1055.1904
$ws.Range("K132").Value = 3365.64
$ws.Range("L132").Value = 3165.5712
$ws.Range("M132").Value = -835.6400000000003
$ws.Range("N132").Value = -8225.5712
$ws.Range("H136").Value = 558.2069
$ws.Range("I136").Value = 515.36365
$ws.Range("J136").Value = 692.8570999999999
$ws.Range("K136").Value = 1546.09095
$ws.Range("L136").Value = 2078.5713
$ws.Range("M136").Value = 1003.90905
$ws.Range("N136").Value = -7178.5713
